$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5; existing rows 5-9 shift down to 6-10.
$ws.Rows.Item(5).Insert()

# Copy formatting of the date cell (D6, which was the old D5) into the new D5
$ws.Range("D6").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 5 with its values
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44799
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112012
$ws.Range("G5").Value = "Espinaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("N5").Value = "$/cuna 10 kilos"
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 700
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = "Hortaliza"
